# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The "Periodo Mora" (arrears period) list in rows 16-28, column E, is
# reordered so the most recent periods (2103, 2102, 2101, ...) come
# first and the oldest (2003) comes last - i.e. the previous
# chronological order is reversed. The "Salario Basico" value (column
# F) travels with its row, so after the reversal the values that used
# to sit on the first/last row (35112 / 26919) end up swapped.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$periods = @("2103", "2102", "2101", "2012", "2011", "2010", "2009", "2008", "2007", "2006", "2005", "2004", "2003")

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
}

# First/last "Salario Basico" values swap as a consequence of the reversal above.
$ws.Range("F16").Value = 26919
$ws.Range("F28").Value = 35112
